# Added filtering options for the Component Analysis
# This clears out specific cells in rows 2,3,5,6,7 (columns I,J,K as applicable)
# to match the updated Component Analysis filtering behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J2:K2").ClearContents()
$ws.Range("I3:K3").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("J6:K6").ClearContents()
$ws.Range("I7:K7").ClearContents()
